$wb = $excel.ActiveWorkbook

# Sheet "Overview" - Latest HO Xliff Generate Date for row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 19:11:26"

# Sheet "zh-cn" - Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 19:11:21"
$wsZhCn.Range("K2").Value = "2016-08-27 19:11:38"

# Sheet "de-de" - Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 19:11:26"
$wsDeDe.Range("K2").Value = "2016-08-27 19:11:45"
